$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data rows (rows 2 and 3), keeping only the header row
$ws.Rows.Item(2).Resize(2).Delete() | Out-Null
Write-Host "rows deleted"

# Shift the existing header cells (A1:E1) one column to the right (B1:F1),
# working from the rightmost column down to the leftmost so nothing gets
# overwritten before it is copied. Each header's formatting (bold header
# style) travels with it.
for ($col = 5; $col -ge 1; $col--) {
    $srcCell = $ws.Cells.Item(1, $col)
    $dstCell = $ws.Cells.Item(1, $col + 1)

    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    Write-Host "formats shifted for col $col"

    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    Write-Host "values shifted for col $col"
}
$excel.CutCopyMode = $false

# Put the new "Unnamed: 0" header into the now-vacated A1 (it already has
# the correct bold header style, since that cell has not moved)
$ws.Range("A1").Value = "Unnamed: 0"
Write-Host "header set"
